$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 259.222432
$ws.Range("H2").Value = 777.667296
$ws.Range("I2").Value = 0.7751755260536969
$ws.Range("J2").Value = 0.7751755260536969
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 18.98606933333333
$ws.Range("N2").Value = 56.958208
$ws.Range("O2").Value = 0.879720792300607
$ws.Range("P2").Value = 0.879720792300607
$ws.Range("Q2").Value = 4921.615066707284
$ws.Range("R2").Value = 44294.53560036557
$ws.Range("S2").Value = 0.681938027951998
$ws.Range("T2").Value = 0.681938027951998

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 259.222432
$ws.Range("H3").Value = 777.667296
$ws.Range("I3").Value = 0.7751755260536969
$ws.Range("J3").Value = 0.7751755260536969
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.122490333333333
$ws.Range("N3").Value = 3.367471
$ws.Range("O3").Value = 0.0520106646643328
$ws.Range("P3").Value = 0.0520106646643328
$ws.Range("Q3").Value = 290.9746741031573
$ws.Range("R3").Value = 2618.772066928416
$ws.Range("S3").Value = 0.0403173943415766
$ws.Range("T3").Value = 0.0403173943415766

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 259.222432
$ws.Range("H4").Value = 777.667296
$ws.Range("I4").Value = 0.7751755260536969
$ws.Range("J4").Value = 0.7751755260536969
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.473366666666666
$ws.Range("N4").Value = 4.4201
$ws.Range("O4").Value = 0.06826854303506025
$ws.Range("P4").Value = 0.06826854303506025
$ws.Range("Q4").Value = 381.9296905610666
$ws.Range("R4").Value = 3437.3672150496
$ws.Range("S4").Value = 0.05292010376012227
$ws.Range("T4").Value = 0.05292010376012227

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 48.69324
$ws.Range("H5").Value = 146.07972
$ws.Range("I5").Value = 0.1456116572976945
$ws.Range("J5").Value = 0.1456116572976945
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 18.98606933333333
$ws.Range("N5").Value = 56.958208
$ws.Range("O5").Value = 0.879720792300607
$ws.Range("P5").Value = 0.879720792300607
$ws.Range("Q5").Value = 924.49323070464
$ws.Range("R5").Value = 8320.43907634176
$ws.Range("S5").Value = 0.1280976025261323
$ws.Range("T5").Value = 0.1280976025261323

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 48.69324
$ws.Range("H6").Value = 146.07972
$ws.Range("I6").Value = 0.1456116572976945
$ws.Range("J6").Value = 0.1456116572976945
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.122490333333333
$ws.Range("N6").Value = 3.367471
$ws.Range("O6").Value = 0.0520106646643328
$ws.Range("P6").Value = 0.0520106646643328
$ws.Range("Q6").Value = 54.65769119868001
$ws.Range("R6").Value = 491.9192207881201
$ws.Range("S6").Value = 0.007573359078928138
$ws.Range("T6").Value = 0.007573359078928137

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 48.69324
$ws.Range("H7").Value = 146.07972
$ws.Range("I7").Value = 0.1456116572976945
$ws.Range("J7").Value = 0.1456116572976945
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.473366666666666
$ws.Range("N7").Value = 4.4201
$ws.Range("O7").Value = 0.06826854303506025
$ws.Range("P7").Value = 0.06826854303506025
$ws.Range("Q7").Value = 71.74299670799999
$ws.Range("R7").Value = 645.686970372
$ws.Range("S7").Value = 0.009940695692634104
$ws.Range("T7").Value = 0.009940695692634102

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 26.48914766666666
$ws.Range("H8").Value = 79.46744299999999
$ws.Range("I8").Value = 0.07921281664860853
$ws.Range("J8").Value = 0.07921281664860852
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 18.98606933333333
$ws.Range("N8").Value = 56.958208
$ws.Range("O8").Value = 0.879720792300607
$ws.Range("P8").Value = 0.879720792300607
$ws.Range("Q8").Value = 502.9247941802382
$ws.Range("R8").Value = 4526.323147622144
$ws.Range("S8").Value = 0.06968516182247661
$ws.Range("T8").Value = 0.0696851618224766

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 26.48914766666666
$ws.Range("H9").Value = 79.46744299999999
$ws.Range("I9").Value = 0.07921281664860853
$ws.Range("J9").Value = 0.07921281664860852
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.122490333333333
$ws.Range("N9").Value = 3.367471
$ws.Range("O9").Value = 0.0520106646643328
$ws.Range("P9").Value = 0.0520106646643328
$ws.Range("Q9").Value = 29.73381219407255
$ws.Range("R9").Value = 267.604309746653
$ws.Range("S9").Value = 0.004119911243828056
$ws.Range("T9").Value = 0.004119911243828055

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 26.48914766666666
$ws.Range("H10").Value = 79.46744299999999
$ws.Range("I10").Value = 0.07921281664860853
$ws.Range("J10").Value = 0.07921281664860852
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.473366666666666
$ws.Range("N10").Value = 4.4201
$ws.Range("O10").Value = 0.06826854303506025
$ws.Range("P10").Value = 0.06826854303506025
$ws.Range("Q10").Value = 39.02822720047777
$ws.Range("R10").Value = 351.2540448043
$ws.Range("S10").Value = 0.005407743582303869
$ws.Range("T10").Value = 0.005407743582303868
